$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("NegativeLoginTest")

# Fix the locator values in NegativeLoginTest sheet: replace stray
# whitespace / mismatched student id strings with the correct,
# consistent values used elsewhere in the sheet.
$ws1.Range("A2").Value = "Student-5"
$ws1.Range("B2").Value = "S1234"

$ws1.Range("A3").Value = "Student-5"
$ws1.Range("B3").Value = "S1234"

$ws1.Range("B4").Value = "S1234"

$ws1.Range("A5").Value = "Student-5"
$ws1.Range("B5").Value = "S1234"
